$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "E7"  = 16.645
    "A9"  = -21.726
    "E12" = 17.646
    "E14" = 17.007
    "A18" = -22.095
    "A20" = -19.915
    "E26" = 16.525
    "A27" = -21.918
    "E27" = 16.531
    "E29" = 16.941
    "A35" = -19.823
    "E37" = 16.855
    "E38" = 16.741
    "E51" = 16.65
    "E52" = 16.657
    "E55" = 16.494
    "A69" = -21.557
    "E69" = 17.438
    "E70" = 17.609
    "A76" = -20.047
    "A78" = -19.854
    "E81" = 16.434
    "A82" = -21.986
    "A83" = -21.891
    "E83" = 16.578
    "A93" = -21.524
    "E102" = 16.724
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
